# This workbook's data rows (2-8) got reshuffled: the values in columns
# D, L, M, N, O, P, Q, R, S, T move between rows while A, B, C, E, F, G, H, I, J, K
# (which are identical across all data rows) stay put. Row 5 is unchanged.
#
# Mapping of new row <- old row (source of the row's D..T values):
#   2 <- 3
#   3 <- 7
#   4 <- 8
#   5 <- 5 (unchanged)
#   6 <- 4
#   7 <- 2
#   8 <- 6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values (columns D and L..T) for data rows 2-8 before
# overwriting anything, since several rows swap values with each other.
$cols = @("D","L","M","N","O","P","Q","R","S","T")

$original = @{}
foreach ($r in 2..8) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $original[$r] = $rowVals
}

# New row -> source row mapping
$mapping = @{
    2 = 3
    3 = 7
    4 = 8
    5 = 5
    6 = 4
    7 = 2
    8 = 6
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $original[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcVals[$col]
    }
}
